$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.250.43"
$ws.Range("E2").Value = "  +1.15%  "
$ws.Range("D3").Value = "2.353.77"
$ws.Range("E3").Value = "  +1.55%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'539.89"
$ws.Range("E5").Value = "  +1.52%  "
$ws.Range("D6").Value = "'135.38"
$ws.Range("E6").Value = "  +2.02%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("E8").Value = "  +4.91%  "
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("E10").Value = "  +5.08%  "
$ws.Range("E11").Value = "  -0.77%  "
$ws.Range("D12").Value = "'0.354"
$ws.Range("E12").Value = "  +1.83%  "
$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D13").Value = "'23.81"
$ws.Range("E13").Value = "  +1.32%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.772.52"
$ws.Range("E14").Value = "  +1.29%  "
$ws.Range("D15").Value = "58.212.35"
$ws.Range("E15").Value = "  +1.37%  "
$ws.Range("D16").Value = "'0.0000133"
$ws.Range("E16").Value = "  +0.42%  "
$ws.Range("D17").Value = "2.346.48"
$ws.Range("E17").Value = "  +0.75%  "
$ws.Range("E18").Value = "  +2.83%  "
$ws.Range("D19").Value = "'332.59"
$ws.Range("E19").Value = "  -1.60%  "
$ws.Range("D20").Value = "'4.28"
$ws.Range("E20").Value = "  +2.59%  "
$ws.Range("D21").Value = "'6.77"
$ws.Range("E21").Value = "  -1.64%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").Value = "'62.88"
$ws.Range("E23").Value = "  +1.10%  "
$ws.Range("E24").Value = "  -0.44%  "
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  +0.75%  "
$ws.Range("E26").Value = "  -2.84%  "
$ws.Range("E27").Value = "  +3.32%  "
$ws.Range("D28").Value = "'172.13"
$ws.Range("E28").Value = "  -0.50%  "
$ws.Range("E29").Value = "  +1.19%  "
$ws.Range("D30").Value = "0.0₃0736"
$ws.Range("E30").Value = "  +1.32%  "
$ws.Range("D31").Value = "'6.13"
$ws.Range("E31").Value = "  +0.15%  "
$ws.Range("D32").Value = "'1.03"
$ws.Range("E32").Value = "  +12.21%  "
$ws.Range("D33").Value = "'18.44"
$ws.Range("E33").Value = "  -0.65%  "
$ws.Range("D35").Value = "'4.26"
$ws.Range("E35").Value = "  +6.29%  "
$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  +0.22%  "
$ws.Range("E37").Value = "  -1.00%  "
$ws.Range("E38").Value = "  +3.62%  "
$ws.Range("D39").Value = "'39.21"
$ws.Range("E39").Value = "  -0.08%  "
$ws.Range("D40").Value = "'145.42"
$ws.Range("E40").Value = "  -2.57%  "
$ws.Range("D41").Value = "'293.98"
$ws.Range("E41").Value = "  +4.52%  "
$ws.Range("D42").Value = "'0.377"
$ws.Range("E42").Value = "  +0.30%  "
$ws.Range("D43").Value = "'3.64"
$ws.Range("E43").Value = "  +0.92%  "
$ws.Range("E44").Value = "  +1.90%  "
$ws.Range("D45").Value = "'19.18"
$ws.Range("E45").Value = "  +1.08%  "
$ws.Range("D46").Value = "'0.0502"
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("D47").Value = "'0.562"
$ws.Range("E47").Value = "  +0.78%  "
$ws.Range("E48").Value = "  +1.23%  "
$ws.Range("D49").Value = "'0.385"
$ws.Range("E49").Value = "  +0.88%  "
$ws.Range("D50").Value = "'17.52"
$ws.Range("E50").Value = "  +0.46%  "
$ws.Range("D51").Value = "'11.07"
$ws.Range("E51").Value = "  +0.49%  "
